$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Model sheet: fill in the newly-disclosed quarterly Shops data
#    (Company-owned row 7, Franchise row 8) for Q322..Q124 and Q424..Q125.
#    Revenue (row 9), Revenue/Shop (row 5) are formula-driven and will
#    recompute automatically.
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("Model")

$model.Range("C7").Value = 174
$model.Range("D7").Value = 176
$model.Range("E7").Value = 173
$model.Range("F7").Value = 221
$model.Range("G7").Value = 236
$model.Range("H7").Value = 227
$model.Range("I7").Value = 248
$model.Range("L7").Value = 314
$model.Range("M7").Value = 326

$model.Range("C8").Value = 25
$model.Range("D8").Value = 26
$model.Range("E8").Value = 24
$model.Range("F8").Value = 29
$model.Range("G8").Value = 28
$model.Range("H8").Value = 27
$model.Range("I8").Value = 27
$model.Range("L8").Value = 28
$model.Range("M8").Value = 29

# ---------------------------------------------------------------------
# 2. Model sheet: new "Shop Growth y/o/y %" row (row 22) mirroring the
#    existing "Revenue Growth y/o/y %" row (row 21), and extend row 21
#    across the newly-populated columns.
# ---------------------------------------------------------------------
$model.Range("C21:O21").NumberFormat = "0.00%"
$model.Range("G21").Formula = "=(G9-C9)/C9"
$model.Range("H21:M21").Formula = "=(H9-D9)/D9"
$model.Range("N21").Formula = "=(N9-J9)/J9"
$model.Range("O21").Formula = "=(O9-K9)/K9"

$model.Range("B22").Value = "Shop Growth y/o/y %"
$model.Range("B22").Style = $model.Range("B21").Style
$model.Range("C22:O22").NumberFormat = "0.00%"
$model.Range("G22").Formula = "=(G4-C4)/C4"
$model.Range("H22:N22").Formula = "=(H4-D4)/D4"
$model.Range("O22").Formula = "=(O4-K4)/K4"

# ---------------------------------------------------------------------
# 3. New "Annuals" sheet, appended after "Model", keep "Model" active.
# ---------------------------------------------------------------------
$annuals = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$annuals.Name = "Annuals"
$model.Activate()

# ---------------------------------------------------------------------
# 4. Main sheet: add an A1 "Model" hyperlink back to the Model tab,
#    mirroring the existing Model!A1 -> Main hyperlink.
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("Main")
$main.Range("A1").Value = "Model"
$main.Hyperlinks.Add($main.Range("A1"), "", "Model!A1", "", "Model")
$main.Range("A1").Style = "Hyperlink"
